$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 7-row block that previously held the "extra" GeeksForGeeks links (rows 23-29)
# is removed entirely; the LeetCode block that followed (old rows 30-35) shifts up
# to rows 23-28, and everything below shifts up to fill the rest (old 36-57 -> 29-50).
$ws.Range("A23:A29").EntireRow.Delete() | Out-Null

# Mark every LeetCode row (now rows 23-28) as completed - "Week 3 completed".
$ws.Range("B23").Value2 = "YES"
$ws.Range("B24").Value2 = "YES"
$ws.Range("B25").Value2 = "YES"
$ws.Range("B26").Value2 = "YES"
$ws.Range("B27").Value2 = "YES"
$ws.Range("B28").Value2 = "YES"

# Reflect where the user's cursor/selection ended up after the edit.
$ws.Range("B27").Select() | Out-Null
